# Fruta / hortaliza, semanal
# Insert 4 weekly report rows (Cara cara + Lane Late, week of 44461) ahead of
# the existing Fukumoto/New Hall/Valencia rows, pushing those down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 286, shifting old rows 286:293 down to 290:297.
$ws.Rows("286:289").Insert()

# Columns that are constant across every row of this block.
$ws.Range("A286:A289").Value = 2
$ws.Range("B286:B289").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C286:C289").Value = "Coquimbo"
$ws.Range("D286:D289").Value = 44461
$ws.Range("E286:E289").Value = 4
$ws.Range("F286:F289").Value = "Fruta"
$ws.Range("G286:G289").Value = 100102
$ws.Range("H286:H289").Value = "Cítricos"
$ws.Range("I286:I289").Value = 100102005
$ws.Range("J286:J289").Value = "Naranja"
$ws.Range("M286:M289").Value = 20
$ws.Range("Q286:Q289").Value = "$/bins (400 kilos)"
$ws.Range("R286:R289").Value = "Provincia de Limarí"
$ws.Range("T286:T289").Value = 400

# Row 286: Cara cara / Primera
$ws.Range("K286").Value = "Cara cara"
$ws.Range("L286").Value = "Primera"
$ws.Range("N286").Value = 110000
$ws.Range("O286").Value = 120000
$ws.Range("P286").Value = 115000
$ws.Range("S286").Value = 288

# Row 287: Cara cara / Segunda
$ws.Range("K287").Value = "Cara cara"
$ws.Range("L287").Value = "Segunda"
$ws.Range("N287").Value = 90000
$ws.Range("O287").Value = 100000
$ws.Range("P287").Value = 95000
$ws.Range("S287").Value = 238

# Row 288: Lane Late / Primera
$ws.Range("K288").Value = "Lane Late"
$ws.Range("L288").Value = "Primera"
$ws.Range("N288").Value = 100000
$ws.Range("O288").Value = 110000
$ws.Range("P288").Value = 105000
$ws.Range("S288").Value = 262

# Row 289: Lane Late / Segunda
$ws.Range("K289").Value = "Lane Late"
$ws.Range("L289").Value = "Segunda"
$ws.Range("N289").Value = 80000
$ws.Range("O289").Value = 90000
$ws.Range("P289").Value = 85000
$ws.Range("S289").Value = 212
